$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing quarter columns (E:I) five columns to the right (-> J:N),
# opening up five fresh columns at E:I for the new, earlier quarters.
$ws.Columns("E:I").Insert()

# Header rows: re-label the freshly opened quarter-header cells.
$ws.Range("E8").Value = "فصل اول منتهی به 1399/03"
$ws.Range("F8").Value = "فصل دوم منتهی به 1399/06"
$ws.Range("G8").Value = "فصل سوم منتهی به 1399/09"
$ws.Range("H8").Value = "فصل چهارم منتهی به 1399/12"
$ws.Range("I8").Value = "فصل اول منتهی به 1400/03"

$ws.Range("E24").Value = "فصل اول منتهی به 1399/03"
$ws.Range("F24").Value = "فصل دوم منتهی به 1399/06"
$ws.Range("G24").Value = "فصل سوم منتهی به 1399/09"
$ws.Range("H24").Value = "فصل چهارم منتهی به 1399/12"
$ws.Range("I24").Value = "فصل اول منتهی به 1400/03"

# Row 10 - هزینه حمل و نقل و انتقال
$ws.Range("E10").Value = 84595
$ws.Range("F10").Value = 260546
$ws.Range("G10").Value = 207394
$ws.Range("H10").Value = 315925
$ws.Range("I10").Value = 84595

# Row 11 - هزینه خدمات پس از فروش
$ws.Range("E11:I11").Value = 0

# Row 12 - حق العمل و کمیسیون فروش
$ws.Range("E12:I12").Value = 0

# Row 13 - هزینه تبلیغات
$ws.Range("E13").Value = 620
$ws.Range("F13").Value = 1240
$ws.Range("G13").Value = 705
$ws.Range("H13").Value = -1283
$ws.Range("I13").Value = 637

# Row 14 - هزینه مواد مصرفی
$ws.Range("E14:I14").Value = 0

# Row 15 - هزینه انرژی (آب، برق، گاز و سوخت)
$ws.Range("E15").Value = 65
$ws.Range("F15").Value = 255
$ws.Range("G15").Value = -139
$ws.Range("H15").Value = 515
$ws.Range("I15").Value = 76

# Row 16 - هزینه استهلاک
$ws.Range("E16").Value = 906
$ws.Range("F16").Value = -215
$ws.Range("G16").Value = 92
$ws.Range("H16").Value = 2354
$ws.Range("I16").Value = 563

# Row 17 - هزینه حقوق و دستمزد
$ws.Range("E17").Value = 16726
$ws.Range("F17").Value = 7764
$ws.Range("G17").Value = 8395
$ws.Range("H17").Value = 27867
$ws.Range("I17").Value = 16726

# Row 18 - هزینه مطالبات مشکوک الوصول
$ws.Range("E18:I18").Value = 0

# Row 19 - سایر هزینه ها
$ws.Range("E19").Value = 798
$ws.Range("F19").Value = 39725
$ws.Range("G19").Value = -11575
$ws.Range("H19").Value = 67242
$ws.Range("I19").Value = 252896

# Row 20 - جمع
$ws.Range("E20").Value = 103710
$ws.Range("F20").Value = 309315
$ws.Range("G20").Value = 204872
$ws.Range("H20").Value = 412620
$ws.Range("I20").Value = 355493

# Row 26 - تعداد پرسنل غیر تولیدی شرکت
$ws.Range("E26").Value = 256
$ws.Range("F26").Value = 226
$ws.Range("G26").Value = 226
$ws.Range("H26").Value = 234
$ws.Range("I26").Value = 234

# Row 27 - تعداد پرسنل تولیدی شرکت
$ws.Range("E27").Value = 50
$ws.Range("F27").Value = 26
$ws.Range("G27").Value = 26
$ws.Range("H27").Value = 24
$ws.Range("I27").Value = 24
